# fix: cambia la tabla sensibilidad
# Updates the sensitivity table values in A2:E6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = 264678.4274872866
    "B2" = 220265.610095929
    "C2" = 184793.0682734134
    "D2" = 156431.1092878721
    "E2" = 133700.69611276

    "A3" = 320876.0600554887
    "B3" = 265652.9572145922
    "C3" = 221434.3448355638
    "D3" = 186039.4470370775
    "E3" = 157679.3479733003

    "A4" = 389725.3899023259
    "B4" = 321474.3733763354
    "C4" = 266612.7400939023
    "D4" = 222587.9169773016
    "E4" = 187272.2952660974

    "A5" = 473479.9394129804
    "B5" = 389725.3899023259
    "C5" = 322062.4270923145
    "D5" = 267558.0992372333
    "E5" = 223726.5962144691

    "A6" = 574562.0712581592
    "B6" = 472596.1001033363
    "C6" = 389725.3899023259
    "D6" = 322640.481462321
    "E6" = 268489.3487276214
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
